$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "increasing room size" sheet ---
$ws3 = $wb.Worksheets.Item("increasing room size")
$ws3.Range("A1:L16").Select()

# --- Add the new "limiting subject conflict" sheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add($null, $lastSheet)
$newWs.Name = "limiting subject conflict"

# Header row
$newWs.Range("A1").Value = "Year"
$newWs.Range("B1").Value = "Classes"
$newWs.Range("C1").Value = "Profs"
$newWs.Range("D1").Value = "Students"
$newWs.Range("E1").Value = "Times"
$newWs.Range("F1").Value = "Rooms"
$newWs.Range("G1").Value = "Time (s)"
$newWs.Range("H1").Value = "Best"
$newWs.Range("I1").Value = "Experimental"
$newWs.Range("J1").Value = "% Optimality"
$newWs.Range("K1").Value = "Diff Exp"
$newWs.Range("L1").Value = "Diff %Opt"

$newWs.Range("A2").Value = 2000
$newWs.Range("B2").Value = 231
$newWs.Range("C2").Value = 164
$newWs.Range("D2").Value = 1112
$newWs.Range("E2").Value = 17
$newWs.Range("F2").Value = 60
$newWs.Range("G2").Value = 0.00797939300537109
$newWs.Range("H2").Value = 3497
$newWs.Range("I2").Value = 3133
$newWs.Range("J2").Formula = "=I2/H2"
$newWs.Range("K2").Formula = "=I2-'no constraints'!I2"
$newWs.Range("L2").Formula = "=J2-'no constraints'!J2"
$newWs.Range("A3").Value = 2001
$newWs.Range("B3").Value = 222
$newWs.Range("C3").Value = 167
$newWs.Range("D3").Value = 1096
$newWs.Range("E3").Value = 17
$newWs.Range("F3").Value = 59
$newWs.Range("G3").Value = 0.00498628616333007
$newWs.Range("H3").Value = 3542
$newWs.Range("I3").Value = 3189
$newWs.Range("J3").Formula = "=I3/H3"
$newWs.Range("K3").Formula = "=I3-'no constraints'!I3"
$newWs.Range("L3").Formula = "=J3-'no constraints'!J3"
$newWs.Range("A4").Value = 2002
$newWs.Range("B4").Value = 239
$newWs.Range("C4").Value = 159
$newWs.Range("D4").Value = 1090
$newWs.Range("E4").Value = 17
$newWs.Range("F4").Value = 61
$newWs.Range("G4").Value = 0.00405073165893554
$newWs.Range("H4").Value = 3579
$newWs.Range("I4").Value = 3112
$newWs.Range("J4").Formula = "=I4/H4"
$newWs.Range("K4").Formula = "=I4-'no constraints'!I4"
$newWs.Range("L4").Formula = "=J4-'no constraints'!J4"
$newWs.Range("A5").Value = 2003
$newWs.Range("B5").Value = 241
$newWs.Range("C5").Value = 151
$newWs.Range("D5").Value = 1104
$newWs.Range("E5").Value = 17
$newWs.Range("F5").Value = 59
$newWs.Range("G5").Value = 0.00800180435180664
$newWs.Range("H5").Value = 3539
$newWs.Range("I5").Value = 3168
$newWs.Range("J5").Formula = "=I5/H5"
$newWs.Range("K5").Formula = "=I5-'no constraints'!I5"
$newWs.Range("L5").Formula = "=J5-'no constraints'!J5"
$newWs.Range("A6").Value = 2004
$newWs.Range("B6").Value = 265
$newWs.Range("C6").Value = 163
$newWs.Range("D6").Value = 1124
$newWs.Range("E6").Value = 17
$newWs.Range("F6").Value = 51
$newWs.Range("G6").Value = 0.00498652458190918
$newWs.Range("H6").Value = 3700
$newWs.Range("I6").Value = 3339
$newWs.Range("J6").Formula = "=I6/H6"
$newWs.Range("K6").Formula = "=I6-'no constraints'!I6"
$newWs.Range("L6").Formula = "=J6-'no constraints'!J6"
$newWs.Range("A7").Value = 2005
$newWs.Range("B7").Value = 255
$newWs.Range("C7").Value = 156
$newWs.Range("D7").Value = 1127
$newWs.Range("E7").Value = 17
$newWs.Range("F7").Value = 52
$newWs.Range("G7").Value = 0.00600171089172363
$newWs.Range("H7").Value = 3680
$newWs.Range("I7").Value = 3304
$newWs.Range("J7").Formula = "=I7/H7"
$newWs.Range("K7").Formula = "=I7-'no constraints'!I7"
$newWs.Range("L7").Formula = "=J7-'no constraints'!J7"
$newWs.Range("A8").Value = 2006
$newWs.Range("B8").Value = 269
$newWs.Range("C8").Value = 169
$newWs.Range("D8").Value = 1167
$newWs.Range("E8").Value = 17
$newWs.Range("F8").Value = 63
$newWs.Range("G8").Value = 0.00800037384033203
$newWs.Range("H8").Value = 3798
$newWs.Range("I8").Value = 3391
$newWs.Range("J8").Formula = "=I8/H8"
$newWs.Range("K8").Formula = "=I8-'no constraints'!I8"
$newWs.Range("L8").Formula = "=J8-'no constraints'!J8"
$newWs.Range("A9").Value = 2007
$newWs.Range("B9").Value = 283
$newWs.Range("C9").Value = 169
$newWs.Range("D9").Value = 1148
$newWs.Range("E9").Value = 17
$newWs.Range("F9").Value = 62
$newWs.Range("G9").Value = 0.00698137283325195
$newWs.Range("H9").Value = 3862
$newWs.Range("I9").Value = 3480
$newWs.Range("J9").Formula = "=I9/H9"
$newWs.Range("K9").Formula = "=I9-'no constraints'!I9"
$newWs.Range("L9").Formula = "=J9-'no constraints'!J9"
$newWs.Range("A10").Value = 2008
$newWs.Range("B10").Value = 284
$newWs.Range("C10").Value = 175
$newWs.Range("D10").Value = 1213
$newWs.Range("E10").Value = 17
$newWs.Range("F10").Value = 63
$newWs.Range("G10").Value = 0.00799942016601562
$newWs.Range("H10").Value = 3794
$newWs.Range("I10").Value = 3410
$newWs.Range("J10").Formula = "=I10/H10"
$newWs.Range("K10").Formula = "=I10-'no constraints'!I10"
$newWs.Range("L10").Formula = "=J10-'no constraints'!J10"
$newWs.Range("A11").Value = 2009
$newWs.Range("B11").Value = 264
$newWs.Range("C11").Value = 164
$newWs.Range("D11").Value = 1352
$newWs.Range("E11").Value = 17
$newWs.Range("F11").Value = 67
$newWs.Range("G11").Value = 0.00799989700317382
$newWs.Range("H11").Value = 4057
$newWs.Range("I11").Value = 3700
$newWs.Range("J11").Formula = "=I11/H11"
$newWs.Range("K11").Formula = "=I11-'no constraints'!I11"
$newWs.Range("L11").Formula = "=J11-'no constraints'!J11"
$newWs.Range("A12").Value = 2010
$newWs.Range("B12").Value = 288
$newWs.Range("C12").Value = 174
$newWs.Range("D12").Value = 1475
$newWs.Range("E12").Value = 17
$newWs.Range("F12").Value = 68
$newWs.Range("G12").Value = 0.00598740577697753
$newWs.Range("H12").Value = 4466
$newWs.Range("I12").Value = 4091
$newWs.Range("J12").Formula = "=I12/H12"
$newWs.Range("K12").Formula = "=I12-'no constraints'!I12"
$newWs.Range("L12").Formula = "=J12-'no constraints'!J12"
$newWs.Range("A13").Value = 2011
$newWs.Range("B13").Value = 280
$newWs.Range("C13").Value = 172
$newWs.Range("D13").Value = 1600
$newWs.Range("E13").Value = 17
$newWs.Range("F13").Value = 64
$newWs.Range("G13").Value = 0.00799989700317382
$newWs.Range("H13").Value = 4671
$newWs.Range("I13").Value = 4312
$newWs.Range("J13").Formula = "=I13/H13"
$newWs.Range("K13").Formula = "=I13-'no constraints'!I13"
$newWs.Range("L13").Formula = "=J13-'no constraints'!J13"
$newWs.Range("A14").Value = 2012
$newWs.Range("B14").Value = 293
$newWs.Range("C14").Value = 175
$newWs.Range("D14").Value = 1659
$newWs.Range("E14").Value = 17
$newWs.Range("F14").Value = 70
$newWs.Range("G14").Value = 0.00800180435180664
$newWs.Range("H14").Value = 4813
$newWs.Range("I14").Value = 4417
$newWs.Range("J14").Formula = "=I14/H14"
$newWs.Range("K14").Formula = "=I14-'no constraints'!I14"
$newWs.Range("L14").Formula = "=J14-'no constraints'!J14"
$newWs.Range("A15").Value = 2013
$newWs.Range("B15").Value = 320
$newWs.Range("C15").Value = 179
$newWs.Range("D15").Value = 1644
$newWs.Range("E15").Value = 17
$newWs.Range("F15").Value = 69
$newWs.Range("G15").Value = 0.0120017528533935
$newWs.Range("H15").Value = 4739
$newWs.Range("I15").Value = 4385
$newWs.Range("J15").Formula = "=I15/H15"
$newWs.Range("K15").Formula = "=I15-'no constraints'!I15"
$newWs.Range("L15").Formula = "=J15-'no constraints'!J15"
$newWs.Range("A16").Value = 2014
$newWs.Range("B16").Value = 280
$newWs.Range("C16").Value = 183
$newWs.Range("D16").Value = 1635
$newWs.Range("E16").Value = 17
$newWs.Range("F16").Value = 67
$newWs.Range("G16").Value = 0.00698137283325195
$newWs.Range("H16").Value = 4558
$newWs.Range("I16").Value = 4165
$newWs.Range("J16").Formula = "=I16/H16"
$newWs.Range("K16").Formula = "=I16-'no constraints'!I16"
$newWs.Range("L16").Formula = "=J16-'no constraints'!J16"

# --- Final selection on the new sheet ---
$newWs.Range("M12").Select()
